$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain formatted as Text so values like "29.382.75"
# or "0.9985" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.382.75"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.847.45"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "240.51"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "0.6306"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.07556"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.2955"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "24.39"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "0.07693"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").Value = "1.849.78"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "4.989"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").Value = "0.6838"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "0.00001000"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "83.01"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "2.106.68"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "6.136"
$ws.Range("E18").Value = "  -2.44%  "
$ws.Range("D19").Value = "29.426.04"
$ws.Range("D20").Value = "227.50"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").Value = "7.548"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "156.92"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").Value = "0.1394"
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").Value = "8.380"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "17.67"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").Value = "1.466"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").Value = "0.05697"
$ws.Range("E30").Value = "  -4.15%  "
$ws.Range("D31").Value = "1.257"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").Value = "4.124"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "4.017"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").Value = "1.846"
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("D35").Value = "1.154"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").Value = "0.7151"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").Value = "2.591"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "1.250.71"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "0.01807"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").Value = "2.781"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").Value = "0.9061"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").Value = "6.179"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "101.20"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "66.10"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").Value = "7.100"
$ws.Range("E47").Value = "  -4.06%  "
$ws.Range("D48").Value = "0.4013"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").Value = "9.070"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").Value = "1.678"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").Value = "0.1119"
